# "Generate Report for Archive"
#
# 1. Update the localization status text "Ready for handoff" -> "In Translation"
#    everywhere it appears:
#      - Overview sheet: E2 (zh-cn) and F2 (de-de)
#      - zh-cn sheet:     C2 (Status column)
#      - de-de sheet:     C2 (Status column)
# 2. Narrow columns E:F on the Overview sheet and column C (Status) on the
#    zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the status text ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- Narrow the now-shorter status columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
